$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Flagged")

# Update Date_Created / Date_Expired timestamps for rows 2-21 on Sheet1
$ws1.Cells.Item(2, 3).Value = 43804.62040354496
$ws1.Cells.Item(2, 4).Value = 44170.62040354496
$ws1.Cells.Item(3, 3).Value = 43804.62042361854
$ws1.Cells.Item(3, 4).Value = 44170.62042361854
$ws1.Cells.Item(4, 3).Value = 43804.62044560831
$ws1.Cells.Item(4, 4).Value = 44170.62044560831
$ws1.Cells.Item(5, 3).Value = 43804.62046317705
$ws1.Cells.Item(5, 4).Value = 44170.62046317705
$ws1.Cells.Item(6, 3).Value = 43804.62048405896
$ws1.Cells.Item(6, 4).Value = 44170.62048405896
$ws1.Cells.Item(7, 3).Value = 43804.62050764143
$ws1.Cells.Item(7, 4).Value = 44170.62050764143
$ws1.Cells.Item(8, 3).Value = 43804.62053384448
$ws1.Cells.Item(8, 4).Value = 44170.62053384448
$ws1.Cells.Item(9, 3).Value = 43804.62055432207
$ws1.Cells.Item(9, 4).Value = 44170.62055432207
$ws1.Cells.Item(10, 3).Value = 43804.62057781933
$ws1.Cells.Item(10, 4).Value = 44170.62057781933
$ws1.Cells.Item(11, 3).Value = 43804.62059922086
$ws1.Cells.Item(11, 4).Value = 44170.62059922086
$ws1.Cells.Item(12, 3).Value = 43804.62062163772
$ws1.Cells.Item(12, 4).Value = 44170.62062163772
$ws1.Cells.Item(13, 3).Value = 43804.62064129587
$ws1.Cells.Item(13, 4).Value = 44170.62064129587
$ws1.Cells.Item(14, 3).Value = 43804.62067017692
$ws1.Cells.Item(14, 4).Value = 44170.62067017692
$ws1.Cells.Item(15, 3).Value = 43804.62068978882
$ws1.Cells.Item(15, 4).Value = 44170.62068978882
$ws1.Cells.Item(16, 3).Value = 43804.62071159387
$ws1.Cells.Item(16, 4).Value = 44170.62071159387
$ws1.Cells.Item(17, 3).Value = 43804.62073176025
$ws1.Cells.Item(17, 4).Value = 44170.62073176025
$ws1.Cells.Item(18, 3).Value = 43804.62076025998
$ws1.Cells.Item(18, 4).Value = 44170.62076025998
$ws1.Cells.Item(19, 3).Value = 43804.6207817073
$ws1.Cells.Item(19, 4).Value = 44170.6207817073
$ws1.Cells.Item(20, 3).Value = 43804.62080901489
$ws1.Cells.Item(20, 4).Value = 44170.62080901489
$ws1.Cells.Item(21, 3).Value = 43804.62082906438
$ws1.Cells.Item(21, 4).Value = 44170.62082906438

# Add the newly-flagged row (row 3) to the Flagged sheet
$ws3.Cells.Item(3, 1).Value = "Edwards"
$ws3.Cells.Item(3, 2).Value = "Keith"
$ws3.Cells.Item(3, 3).Value = "2019-12-05 14:53:53.686462"
$ws3.Cells.Item(3, 4).Value = "2020-12-05 14:53:53.686462"
$ws3.Cells.Item(3, 5).Value = "MS Center of NE New York"
$ws3.Cells.Item(3, 6).Value = "Latham, NY"
